# Update the NATMI ligand-receptor TPM-derived metrics on Sheet1
# (rows 2-4) to the newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (ECs -> Efnb2/Rhbdl2 -> ECs)
$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("M2").Value = 0.2203263333333333
$ws.Range("N2").Value = 0.660979
$ws.Range("Q2").Value = 9.41956957322822
$ws.Range("R2").Value = 84.77612615905399
$ws.Range("S2").Value = 0.8529286054750734
$ws.Range("T2").Value = 0.8529286054750735

# Row 3 (FAPs -> Efnb2/Rhbdl2 -> ECs)
$ws.Range("I3").Value = 0.04642608686423023
$ws.Range("J3").Value = 0.04642608686423023
$ws.Range("M3").Value = 0.2203263333333333
$ws.Range("N3").Value = 0.660979
$ws.Range("Q3").Value = 0.5127202352262222
$ws.Range("R3").Value = 4.614482117036
$ws.Range("S3").Value = 0.04642608686423023
$ws.Range("T3").Value = 0.04642608686423023

# Row 4 (MuSCs -> Efnb2/Rhbdl2 -> ECs)
$ws.Range("G4").Value = 5.044817999999999
$ws.Range("I4").Value = 0.1006453076606963
$ws.Range("J4").Value = 0.1006453076606963
$ws.Range("M4").Value = 0.2203263333333333
$ws.Range("N4").Value = 0.660979
$ws.Range("Q4").Value = 1.111506252274
$ws.Range("R4").Value = 10.003556270466
$ws.Range("S4").Value = 0.1006453076606963
$ws.Range("T4").Value = 0.1006453076606963
